$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the newly trained model's (row 4, MobileNetV2(alpha=0.35)) Accuracy/Loss values
$ws.Range("E4").Value = 0.4514
$ws.Range("F4").Value = 0.8462

# Fill in row 5 (MobileNetV2(alpha=0.50)) Accuracy/Loss values
$ws.Range("E5").Value = 0.457
$ws.Range("F5").Value = 0.8279

# Update the active selection to F5, matching the author's last edit location
$ws.Range("F5").Select()
